# "change size of test case"
# Resizes the scaling-test grid on Sheet1: shrinks the big per-cell extents
# (K/L/M/N/O/P columns) by a factor of 10 on X/Y and a factor of 4 on Z,
# turns B2/C2/D2 into formulas derived from those extents, and adds a new
# "UNIT" label column (I2:I6) describing each row's grid resolution.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 -----------------------------------------------------------
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 1600
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = 400
$ws.Range("P2").Value = 0.5

$ws.Range("B2").Formula = "=K2*L2"
$ws.Range("C2").Formula = "=M2*N2"
$ws.Range("D2").Formula = "=O2*P2"

$ws.Range("I2").Value = "512n"
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").VerticalAlignment = -4108

# --- Row 3 -----------------------------------------------------------
$ws.Range("K3").Value = 800
$ws.Range("L3").Value = 20
$ws.Range("M3").Value = 800
$ws.Range("N3").Value = 20
$ws.Range("O3").Value = 200
$ws.Range("P3").Value = 1

$ws.Range("I3").Value = "64n"
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I3").VerticalAlignment = -4108

# --- Row 4 -----------------------------------------------------------
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 40
$ws.Range("M4").Value = 400
$ws.Range("N4").Value = 40
$ws.Range("O4").Value = 100
$ws.Range("P4").Value = 2

$ws.Range("I4").Value = "8n"
$ws.Range("I4").HorizontalAlignment = -4108
$ws.Range("I4").VerticalAlignment = -4108

# --- Row 5 -----------------------------------------------------------
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 80
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 80
$ws.Range("O5").Value = 50
$ws.Range("P5").Value = 4

$ws.Range("I5").Value = "1n"
$ws.Range("I5").HorizontalAlignment = -4108
$ws.Range("I5").VerticalAlignment = -4108

# --- Row 6 -----------------------------------------------------------
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 160
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 160
$ws.Range("O6").Value = 25
$ws.Range("P6").Value = 8

$ws.Range("I6").Value = "7p"
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I6").VerticalAlignment = -4108

# --- column J is now wider to fit the new labels ----------------------
$ws.Columns(10).ColumnWidth = 12.4285714285714

# --- selection moves to the newly-edited cell --------------------------
$ws.Range("I6").Select()
